# Generate Report for Handoff
# Updates the localization-status report after a new handoff xliff was generated
# for the "1b99fc67-4e3b-4c71-90cc-a813252568ae" file (priority bumped from
# "low" to "ht", and the handoff timestamps refreshed).

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet -----------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("H4").Value = "2016-09-05 06:34:46"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("H5").Value = "2016-09-05 06:34:46"
$wsZh.Range("E6").Value = "ht"
$wsZh.Range("H6").Value = "2016-09-05 06:34:46"
$wsZh.Range("E7").Value = "ht"
$wsZh.Range("H7").Value = "2016-09-05 06:34:46"

# --- de-de sheet -------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("H4").Value = "2016-09-05 06:34:52"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("H5").Value = "2016-09-05 06:34:52"
$wsDe.Range("E6").Value = "ht"
$wsDe.Range("H6").Value = "2016-09-05 06:34:52"
$wsDe.Range("E7").Value = "ht"
$wsDe.Range("H7").Value = "2016-09-05 06:34:52"

# --- Overview sheet (shares the same "Latest HO Xliff Generate Date" string) -
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-05 06:34:52"
$wsOverview.Range("G5").Value = "2016-09-05 06:34:52"
$wsOverview.Range("G6").Value = "2016-09-05 06:34:52"
$wsOverview.Range("G7").Value = "2016-09-05 06:34:52"
